$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the text of cell D5 (shared string describing the ensemble architecture)
$ws.Range("D5").Value = "An ensemble of 5 Neural Networks each having the architecture:
1. Embedding Layer
2. BiLSTM Layer
3. GlobalMaxPooling1D 
4. BatchNormalization
5. Drop out of 0.5
6. Dense Layer with relu activation
7. Drop out of 0.5
8. Dense with with relu activation
9. Dropout of 0.5
10. Dense with sigmoid activation"

# Update row 5 height (switch from explicit custom height to auto-fit height 153)
$ws.Rows("5").RowHeight = 153

# Update view: scroll so row 4 is the top-left visible row, and change selection to B6
$ws.Range("B6").Select()
$excel.ActiveWindow.ScrollRow = 4
